$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 464.2
$ws.Range("I33").Value = 283.2
$ws.Range("K33").Value = 283.2
$ws.Range("M33").Value = -54.19999999999999
$ws.Range("H88").Value = 1110.5555
$ws.Range("J88").Value = 1199.5
$ws.Range("L88").Value = 1199.5
$ws.Range("N88").Value = -2011.5
$ws.Range("H91").Value = 1110.5555
$ws.Range("J91").Value = 1199.5
$ws.Range("L91").Value = 1199.5
$ws.Range("N91").Value = -4007.5
$ws.Range("H107").Value = 279.84616
$ws.Range("I107").Value = 342.7
$ws.Range("J107").Value = 70.333336
$ws.Range("K107").Value = 342.7
$ws.Range("L107").Value = 70.333336
$ws.Range("M107").Value = 1577.3
$ws.Range("N107").Value = -3910.333336
$ws.Range("H135").Value = 66669830
$ws.Range("I135").Value = 76926480
$ws.Range("J135").Value = 1597.5
$ws.Range("K135").Value = 692338320
$ws.Range("L135").Value = 14377.5
$ws.Range("M135").Value = -692335785
$ws.Range("N135").Value = -19447.5
$ws.Range("H137").Value = 2025.4231
$ws.Range("I137").Value = 1592.2354
$ws.Range("K137").Value = 4776.706200000001
$ws.Range("M137").Value = -2226.706200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1207.1578
$ws.Range("I2").Value = 1242.6666
$ws.Range("K2").Value = 1242.6666
$ws.Range("M2").Value = -1129.6666
$ws.Range("H32").Value = 2617.762
$ws.Range("I32").Value = 1547.5193
$ws.Range("J32").Value = 7677.091
$ws.Range("K32").Value = 1547.5193
$ws.Range("L32").Value = 7677.091
$ws.Range("M32").Value = -1260.5193
$ws.Range("N32").Value = -8251.091
$ws.Range("H45").Value = 1722.3077
$ws.Range("I45").Value = 1349.5
$ws.Range("J45").Value = 2965
$ws.Range("K45").Value = 1349.5
$ws.Range("L45").Value = 2965
$ws.Range("M45").Value = -972.5
$ws.Range("N45").Value = -3719
$ws.Range("H88").Value = 15153689
$ws.Range("J88").Value = 2858.2856
$ws.Range("L88").Value = 2858.2856
$ws.Range("N88").Value = -3670.2856
$ws.Range("H91").Value = 15153689
$ws.Range("J91").Value = 2858.2856
$ws.Range("L91").Value = 2858.2856
$ws.Range("N91").Value = -5666.2856
$ws.Range("H97").Value = 546.7143
$ws.Range("I97").Value = 534.05
$ws.Range("K97").Value = 534.05
$ws.Range("M97").Value = -38.04999999999995
$ws.Range("H116").Value = 1207.1578
$ws.Range("I116").Value = 1242.6666
$ws.Range("K116").Value = 1242.6666
$ws.Range("M116").Value = 1051.3334
$ws.Range("H132").Value = 20835002
$ws.Range("I132").Value = 25001548
$ws.Range("K132").Value = 75004644
$ws.Range("M132").Value = -75002114

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1207.1578
$ws.Range("I3").Value = 1242.6666
$ws.Range("K3").Value = 1242.6666
$ws.Range("M3").Value = -1128.6666
$ws.Range("H81").Value = 81954.336
$ws.Range("J81").Value = 81954.336
$ws.Range("L81").Value = 81954.336
$ws.Range("N81").Value = -84076.336
$ws.Range("H84").Value = 81954.336
$ws.Range("J84").Value = 81954.336
$ws.Range("L84").Value = 245863.008
$ws.Range("N84").Value = -256471.008
$ws.Range("H105").Value = 1661.6316
$ws.Range("I105").Value = 1552.0667
$ws.Range("K105").Value = 1552.0667
$ws.Range("M105").Value = 194.9332999999999
$ws.Range("H128").Value = 8450
$ws.Range("I128").Value = 8450
$ws.Range("K128").Value = 25350
$ws.Range("M128").Value = -22860
$ws.Range("H134").Value = 2440.868
$ws.Range("I134").Value = 2344.4666
$ws.Range("J134").Value = 2983.125
$ws.Range("K134").Value = 7033.399800000001
$ws.Range("L134").Value = 8949.375
$ws.Range("M134").Value = -4498.399800000001
$ws.Range("N134").Value = -14019.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2969.2395
$ws.Range("I31").Value = 1600.8695
$ws.Range("J31").Value = 3624.9167
$ws.Range("K31").Value = 1600.8695
$ws.Range("L31").Value = 3624.9167
$ws.Range("M31").Value = -1305.8695
$ws.Range("N31").Value = -4214.9167
$ws.Range("H34").Value = 2969.2395
$ws.Range("I34").Value = 1600.8695
$ws.Range("J34").Value = 3624.9167
$ws.Range("K34").Value = 1600.8695
$ws.Range("L34").Value = 3624.9167
$ws.Range("M34").Value = -1398.8695
$ws.Range("N34").Value = -4028.9167
$ws.Range("H107").Value = 2028.3823
$ws.Range("I107").Value = 1613.0741
$ws.Range("K107").Value = 1613.0741
$ws.Range("M107").Value = 306.9259
$ws.Range("H122").Value = 2683.0833
$ws.Range("J122").Value = 2448.5
$ws.Range("L122").Value = 7345.5
$ws.Range("N122").Value = -12245.5
$ws.Range("H132").Value = 2650.2188
$ws.Range("I132").Value = 2475.2856
$ws.Range("J132").Value = 3874.75
$ws.Range("K132").Value = 7425.8568
$ws.Range("L132").Value = 11624.25
$ws.Range("M132").Value = -4895.8568
$ws.Range("N132").Value = -16684.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 325.58334
$ws.Range("J2").Value = 223.5
$ws.Range("L2").Value = 1341
$ws.Range("N2").Value = -1567
$ws.Range("H138").Value = 3363311
$ws.Range("I138").Value = 4617431
$ws.Range("K138").Value = 13852293
$ws.Range("M138").Value = -13847153

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 6833.3335
$ws.Range("I99").Value = 1250.25
$ws.Range("K99").Value = 1250.25
$ws.Range("M99").Value = 995.75
$ws.Range("H113").Value = 2037.2727
$ws.Range("I113").Value = 1235
$ws.Range("K113").Value = 1235
$ws.Range("M113").Value = 935

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1167
$ws.Range("I22").Value = 1250
$ws.Range("K22").Value = 1250
$ws.Range("M22").Value = -955
$ws.Range("H27").Value = 1167
$ws.Range("I27").Value = 1250
$ws.Range("K27").Value = 1250
$ws.Range("M27").Value = -1143
$ws.Range("H55").Value = 355.9565
$ws.Range("I55").Value = 334.58823
$ws.Range("J55").Value = 416.5
$ws.Range("K55").Value = 334.58823
$ws.Range("L55").Value = 416.5
$ws.Range("M55").Value = -161.58823
$ws.Range("N55").Value = -762.5
$ws.Range("H82").Value = 1262.6
$ws.Range("I82").Value = 897.3333
$ws.Range("K82").Value = 897.3333
$ws.Range("M82").Value = -536.3333
$ws.Range("H85").Value = 1262.6
$ws.Range("I85").Value = 897.3333
$ws.Range("K85").Value = 897.3333
$ws.Range("M85").Value = 350.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 578.52
$ws.Range("I107").Value = 504.25
$ws.Range("K107").Value = 1512.75
$ws.Range("M107").Value = 407.25
$ws.Range("H122").Value = 2532.3333
$ws.Range("I122").Value = 2499.6365
$ws.Range("J122").Value = 2622.25
$ws.Range("K122").Value = 7498.9095
$ws.Range("L122").Value = 7866.75
$ws.Range("M122").Value = -5048.9095
$ws.Range("N122").Value = -12766.75
$ws.Range("H132").Value = 3361.5945
$ws.Range("I132").Value = 3362
$ws.Range("J132").Value = 3360.5
$ws.Range("K132").Value = 10086
$ws.Range("L132").Value = 10081.5
$ws.Range("M132").Value = -7556
$ws.Range("N132").Value = -15141.5
